$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "(A)" rows (Brook Accomando (A), Gabriella Marquez (A)) are gone and the whole
# staff list was reloaded with refreshed numbers, which nets out to two fewer rows overall
# (20 -> 18). Remove the two now-unused trailing rows (19-20) so the row-level formatting
# for rows 1-13 (ht=20) / 14-18 (default) stays aligned with the new, shorter table, then
# overwrite rows 5-18 in place with the reloaded data (this also retires the "(A)" rows,
# since their old positions get overwritten by the next employee's data).
$ws.Range("A19:O20").EntireRow.Delete()

# Reload the remaining staff rows (5-18) with refreshed performance figures.
$ws.Cells.Item(5, 1).Value = 'Brook Accomando'
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 82
$ws.Cells.Item(5, 7).Value = 88.34
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 82
$ws.Cells.Item(5, 13).Value = 88.34
$ws.Cells.Item(5, 14).Value = 41
$ws.Cells.Item(5, 15).Value = 44.17

$ws.Cells.Item(6, 1).Value = 'Chrissy Cummings'
$ws.Cells.Item(6, 2).Value = 27
$ws.Cells.Item(6, 3).Value = 11
$ws.Cells.Item(6, 4).Value = 10
$ws.Cells.Item(6, 5).Value = 2.5
$ws.Cells.Item(6, 6).Value = 1939.82
$ws.Cells.Item(6, 7).Value = 2090.52
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 1939.82
$ws.Cells.Item(6, 13).Value = 2090.52
$ws.Cells.Item(6, 14).Value = 71.84518518518519
$ws.Cells.Item(6, 15).Value = 77.42666666666666

$ws.Cells.Item(7, 1).Value = 'Danielle Mai'
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 2.1
$ws.Cells.Item(7, 6).Value = 491
$ws.Cells.Item(7, 7).Value = 528.73
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 5
$ws.Cells.Item(7, 11).Value = 5.38
$ws.Cells.Item(7, 12).Value = 496
$ws.Cells.Item(7, 13).Value = 534.11
$ws.Cells.Item(7, 14).Value = 82.66666666666667
$ws.Cells.Item(7, 15).Value = 89.01833333333333

$ws.Cells.Item(8, 1).Value = 'Gabriella Marquez'
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 177.3
$ws.Cells.Item(8, 7).Value = 190.97
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 177.3
$ws.Cells.Item(8, 13).Value = 190.97
$ws.Cells.Item(8, 14).Value = 59.1
$ws.Cells.Item(8, 15).Value = 63.65666666666667

$ws.Cells.Item(9, 1).Value = 'Izzy Kruis'
$ws.Cells.Item(9, 2).Value = 36
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 4
$ws.Cells.Item(9, 5).Value = 2.4
$ws.Cells.Item(9, 6).Value = 2493
$ws.Cells.Item(9, 7).Value = 2686.23
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 2493
$ws.Cells.Item(9, 13).Value = 2686.23
$ws.Cells.Item(9, 14).Value = 69.25
$ws.Cells.Item(9, 15).Value = 74.6175

$ws.Cells.Item(10, 1).Value = 'Jasmine Gomez'
$ws.Cells.Item(10, 2).Value = 19
$ws.Cells.Item(10, 3).Value = 18
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 1.9
$ws.Cells.Item(10, 6).Value = 1215.6
$ws.Cells.Item(10, 7).Value = 1309.02
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 1215.6
$ws.Cells.Item(10, 13).Value = 1309.02
$ws.Cells.Item(10, 14).Value = 63.97894736842105
$ws.Cells.Item(10, 15).Value = 68.8957894736842

$ws.Cells.Item(11, 1).Value = 'Jasmine Saiz'
$ws.Cells.Item(11, 2).Value = 48
$ws.Cells.Item(11, 3).Value = 10
$ws.Cells.Item(11, 4).Value = 13
$ws.Cells.Item(11, 5).Value = 2.45
$ws.Cells.Item(11, 6).Value = 3598
$ws.Cells.Item(11, 7).Value = 3876.94
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 3598
$ws.Cells.Item(11, 13).Value = 3876.94
$ws.Cells.Item(11, 14).Value = 74.95833333333333
$ws.Cells.Item(11, 15).Value = 80.76958333333333

$ws.Cells.Item(12, 1).Value = 'Justyne Martinez '
$ws.Cells.Item(12, 2).Value = 24
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 7
$ws.Cells.Item(12, 5).Value = 2.3
$ws.Cells.Item(12, 6).Value = 1833
$ws.Cells.Item(12, 7).Value = 1973.87
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 1833
$ws.Cells.Item(12, 13).Value = 1973.87
$ws.Cells.Item(12, 14).Value = 76.375
$ws.Cells.Item(12, 15).Value = 82.24458333333332

$ws.Cells.Item(13, 1).Value = 'Karen Trevizo'
$ws.Cells.Item(13, 2).Value = 6
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 2.1
$ws.Cells.Item(13, 6).Value = 367
$ws.Cells.Item(13, 7).Value = 395.42
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 367
$ws.Cells.Item(13, 13).Value = 395.42
$ws.Cells.Item(13, 14).Value = 61.16666666666666
$ws.Cells.Item(13, 15).Value = 65.90333333333334

$ws.Cells.Item(14, 1).Value = 'Krisdee Martinez'
$ws.Cells.Item(14, 2).Value = 34
$ws.Cells.Item(14, 3).Value = 12
$ws.Cells.Item(14, 4).Value = 8
$ws.Cells.Item(14, 5).Value = 2.4
$ws.Cells.Item(14, 6).Value = 2270
$ws.Cells.Item(14, 7).Value = 2445.98
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 2270
$ws.Cells.Item(14, 13).Value = 2445.98
$ws.Cells.Item(14, 14).Value = 66.76470588235294
$ws.Cells.Item(14, 15).Value = 71.94058823529411

$ws.Cells.Item(15, 1).Value = 'Maggie  Farrell'
$ws.Cells.Item(15, 2).Value = 31
$ws.Cells.Item(15, 3).Value = 11
$ws.Cells.Item(15, 4).Value = 7
$ws.Cells.Item(15, 5).Value = 2.5
$ws.Cells.Item(15, 6).Value = 2170.8
$ws.Cells.Item(15, 7).Value = 2337.67
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 45.5
$ws.Cells.Item(15, 11).Value = 49
$ws.Cells.Item(15, 12).Value = 2216.3
$ws.Cells.Item(15, 13).Value = 2386.67
$ws.Cells.Item(15, 14).Value = 71.49354838709678
$ws.Cells.Item(15, 15).Value = 76.98935483870967

$ws.Cells.Item(16, 1).Value = 'Makayla Baca'
$ws.Cells.Item(16, 2).Value = 44
$ws.Cells.Item(16, 3).Value = 7
$ws.Cells.Item(16, 4).Value = 14
$ws.Cells.Item(16, 5).Value = 2.5
$ws.Cells.Item(16, 6).Value = 3581
$ws.Cells.Item(16, 7).Value = 3856.23
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 46.5
$ws.Cells.Item(16, 11).Value = 50.07
$ws.Cells.Item(16, 12).Value = 3627.5
$ws.Cells.Item(16, 13).Value = 3906.3
$ws.Cells.Item(16, 14).Value = 82.44318181818181
$ws.Cells.Item(16, 15).Value = 88.77954545454546

$ws.Cells.Item(17, 1).Value = 'Matthew Young'
$ws.Cells.Item(17, 2).Value = 17
$ws.Cells.Item(17, 3).Value = 16
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 2.1
$ws.Cells.Item(17, 6).Value = 1139
$ws.Cells.Item(17, 7).Value = 1227.22
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1139
$ws.Cells.Item(17, 13).Value = 1227.22
$ws.Cells.Item(17, 14).Value = 67
$ws.Cells.Item(17, 15).Value = 72.18941176470588

$ws.Cells.Item(18, 1).Value = 'Vy Torino'
$ws.Cells.Item(18, 2).Value = 37
$ws.Cells.Item(18, 3).Value = 11
$ws.Cells.Item(18, 4).Value = 17
$ws.Cells.Item(18, 5).Value = 2.4
$ws.Cells.Item(18, 6).Value = 2807.5
$ws.Cells.Item(18, 7).Value = 3023.3
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 25
$ws.Cells.Item(18, 11).Value = 26.92
$ws.Cells.Item(18, 12).Value = 2832.5
$ws.Cells.Item(18, 13).Value = 3050.22
$ws.Cells.Item(18, 14).Value = 76.55405405405405
$ws.Cells.Item(18, 15).Value = 82.43837837837837
